# "updated testcases- API-WOrk order"
# The "RMA Details Maintenance Grid" sheet holds generated RMA testcase
# values. Rows 2-4 get re-pointed from the old "RMA-4Z9C" testcase group
# to a fresh "RMA-PPLQ" group (RMA number, line RMA number, and the
# Salesforce record id that the line RMA corresponds to).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

$ws.Range("E2").Value = "RMA-PPLQ-001"
$ws.Range("F2").Value = "RMA-PPLQ-1-1"
$ws.Range("J2").Value = "a7s5f000000xL3IAAU"

$ws.Range("E3").Value = "RMA-PPLQ-002"
$ws.Range("F3").Value = "RMA-PPLQ-1-2"
$ws.Range("J3").Value = "a7s5f000000xL3JAAU"

$ws.Range("E4").Value = "RMA-PPLQ-003"
$ws.Range("F4").Value = "RMA-PPLQ-1-3"
$ws.Range("J4").Value = "a7s5f000000xL3KAAU"
